$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of metric data at the end of the sheet (row 28)
$row = 28

$ws.Range("A$row").Value = "2025-04-28 22:47:54"
$ws.Range("B$row").Value = 334
